$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'242.54"
$ws.Range("E2").Value = "'1BNBBNBBestin24h"

$ws.Range("D4").Value = "'5.524"

$ws.Range("D5").Value = "'0.05588"

$ws.Range("D6").Value = "'3.382"

$ws.Range("D7").Value = "'6.469"

$ws.Range("D8").Value = "'1.075"

$ws.Range("D9").Value = "'0.8028"

$ws.Range("D10").Value = "'0.1417"

$ws.Range("D11").Value = "'0.07448"

$ws.Range("B13").Value = "'ProBitToken"
$ws.Range("C13").Value = "'https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D13").Value = "'0.1305"
$ws.Range("E13").Value = "'12ProBitTokenPROB"

$ws.Range("B14").Value = "'BitrueCoin"
$ws.Range("C14").Value = "'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D14").Value = "'0.02993"
$ws.Range("E14").Value = "'13BitrueCoinBTR"

$ws.Range("B15").Value = "'BitMartToken"
$ws.Range("C15").Value = "'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").Value = "'0.09253"
$ws.Range("E15").Value = "'14BitMartTokenBMX"

$ws.Range("B16").Value = "'BitForexToken"
$ws.Range("C16").Value = "'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").Value = "'0.001662"
$ws.Range("E16").Value = "'15BitForexTokenBF"

$ws.Range("B17").Value = "'MCDex"
$ws.Range("C17").Value = "'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D17").Value = "'3.259"
$ws.Range("E17").Value = "'16MCDexMCB"

$ws.Range("B18").Value = "'CoinExToken"
$ws.Range("C18").Value = "'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D18").Value = "'0.04707"
$ws.Range("E18").Value = "'17CoinExTokenCET"

$ws.Range("B19").Value = "'One"
$ws.Range("C19").Value = "'https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D19").Value = "'0.0005750"
$ws.Range("E19").Value = "'18OneONE"

$ws.Range("B20").Value = "'TigerCash"
$ws.Range("C20").Value = "'https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D20").Value = "'0.006263"
$ws.Range("E20").Value = "'19TigerCashTCH"

$ws.Range("B21").Value = "'BitKan"
$ws.Range("C21").Value = "'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D21").Value = "'0.001052"
$ws.Range("E21").Value = "'20BitKanKAN"

$ws.Range("B22").Value = "'HotbitToken"
$ws.Range("C22").Value = "'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D22").Value = "'0.003803"
$ws.Range("E22").Value = "'21HotbitTokenHTB"

$ws.Range("B23").Value = "'NitroEx"
$ws.Range("C23").Value = "'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D23").Value = "'0.0001500"
$ws.Range("E23").Value = "'22NitroExNTX"

$ws.Range("B24").Value = "'UpBots"
$ws.Range("C24").Value = "'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("D24").Value = "'0.0004777"
$ws.Range("E24").Value = "'23UpBotsUBXT"

$ws.Range("B25").Value = "'LEO"
$ws.Range("C25").Value = "'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D25").Value = "'3.979"
$ws.Range("E25").Value = "'24LEOLEO"

$ws.Range("B26").Value = "'BTSEToken"
$ws.Range("C26").Value = "'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D26").Value = "'2.133"
$ws.Range("E26").Value = "'25BTSETokenBTSE"

$ws.Range("B27").Value = "'BitpandaEcosystemToken"
$ws.Range("C27").Value = "'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D27").Value = "'0.3311"
$ws.Range("E27").Value = "'26BitpandaEcosystemTokenBEST"

$ws.Range("D40").Value = "'0.04176"

$ws.Range("D41").Value = "'0.007009"

$ws.Range("B42").Value = "'BKEXToken"
$ws.Range("C42").Value = "'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "'0.1044"
$ws.Range("E42").Value = "'41BKEXTokenBKK"

$ws.Range("B43").Value = "'CEJI"
$ws.Range("C43").Value = "'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").Value = "'0.002970"
$ws.Range("E43").Value = "'42CEJICEJI"

$ws.Range("D44").Value = "'0.009034"

$ws.Range("D45").Value = "'0.00005492"

$ws.Range("D48").Value = "'0.03025"
$ws.Range("E48").Value = "'47BOLOBOLOWorstin24h"
